$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (28 and 29) so the table shrinks from 28 to 26 data rows
$ws.Rows(28).Delete()
$ws.Rows(28).Delete()

# Rewrite the client / date / amount columns (B, C, D) for rows 2-27 to match the refreshed data
$ws.Cells.Item(2, 2).Value = "Aliso"
$ws.Cells.Item(2, 3).Value = 45983
$ws.Cells.Item(2, 4).Value = 101000

$ws.Cells.Item(3, 2).Value = "Campo Verde Zipaquira"
$ws.Cells.Item(3, 3).Value = 45988
$ws.Cells.Item(3, 4).Value = 64200

$ws.Cells.Item(4, 2).Value = "Carnes Johana"
$ws.Cells.Item(4, 3).Value = 45993
$ws.Cells.Item(4, 4).Value = 176000

$ws.Cells.Item(5, 2).Value = "Carnilandia "
$ws.Cells.Item(5, 3).Value = 45990
$ws.Cells.Item(5, 4).Value = 436900

$ws.Cells.Item(6, 2).Value = "Carnivoros"
$ws.Cells.Item(6, 3).Value = 45959
$ws.Cells.Item(6, 4).Value = 437000

$ws.Cells.Item(7, 2).Value = "Cimarron Dorado"
$ws.Cells.Item(7, 3).Value = 45992
$ws.Cells.Item(7, 4).Value = 407000

$ws.Cells.Item(8, 2).Value = "Cocina China"
$ws.Cells.Item(8, 3).Value = 45992
$ws.Cells.Item(8, 4).Value = 170000

$ws.Cells.Item(9, 2).Value = "Cocina China "
$ws.Cells.Item(9, 3).Value = 45994
$ws.Cells.Item(9, 4).Value = 85000

$ws.Cells.Item(10, 2).Value = "Darwin Futbol"
$ws.Cells.Item(10, 3).Value = 45921
$ws.Cells.Item(10, 4).Value = 200000

$ws.Cells.Item(11, 2).Value = "Davidcito"
$ws.Cells.Item(11, 3).Value = 45947
$ws.Cells.Item(11, 4).Value = 100000

$ws.Cells.Item(12, 2).Value = "El Ruby"
$ws.Cells.Item(12, 3).Value = 45992
$ws.Cells.Item(12, 4).Value = 85100

$ws.Cells.Item(13, 2).Value = "Jordan"
$ws.Cells.Item(13, 3).Value = 45987
$ws.Cells.Item(13, 4).Value = 1313000

$ws.Cells.Item(14, 2).Value = "La Selecta "
$ws.Cells.Item(14, 3).Value = 45912
$ws.Cells.Item(14, 4).Value = 82000

$ws.Cells.Item(15, 2).Value = "Mariana"
$ws.Cells.Item(15, 3).Value = 45650
$ws.Cells.Item(15, 4).Value = 171900

$ws.Cells.Item(16, 2).Value = "Merka Fruver Alejandro"
$ws.Cells.Item(16, 3).Value = 45988
$ws.Cells.Item(16, 4).Value = 60900

$ws.Cells.Item(17, 2).Value = "Merka Fruver Mildred"
$ws.Cells.Item(17, 3).Value = 45988
$ws.Cells.Item(17, 4).Value = 115400

$ws.Cells.Item(18, 2).Value = "Meza 2"
$ws.Cells.Item(18, 3).Value = 45989
$ws.Cells.Item(18, 4).Value = 188000

$ws.Cells.Item(19, 2).Value = "Multicarnes"
$ws.Cells.Item(19, 3).Value = 45989
$ws.Cells.Item(19, 4).Value = 558300

$ws.Cells.Item(20, 2).Value = "Novillon San Mateo"
$ws.Cells.Item(20, 3).Value = 45971
$ws.Cells.Item(20, 4).Value = 83000

$ws.Cells.Item(21, 2).Value = "Pinilla"
$ws.Cells.Item(21, 3).Value = 45931
$ws.Cells.Item(21, 4).Value = 166000

$ws.Cells.Item(22, 2).Value = "Pinilla"
$ws.Cells.Item(22, 3).Value = 45924
$ws.Cells.Item(22, 4).Value = 16000

$ws.Cells.Item(23, 2).Value = "Pinilla Soacha"
$ws.Cells.Item(23, 3).Value = 45993
$ws.Cells.Item(23, 4).Value = 129000

$ws.Cells.Item(24, 2).Value = "Plaza Jessica"
$ws.Cells.Item(24, 3).Value = 45993
$ws.Cells.Item(24, 4).Value = 621000

$ws.Cells.Item(25, 2).Value = "Santander Norte"
$ws.Cells.Item(25, 3).Value = 45973
$ws.Cells.Item(25, 4).Value = 216400

$ws.Cells.Item(26, 2).Value = "Santander Sur"
$ws.Cells.Item(26, 3).Value = 45993
$ws.Cells.Item(26, 4).Value = 80000

$ws.Cells.Item(27, 2).Value = "Vnzlno Punta Anca"
$ws.Cells.Item(27, 3).Value = 45992
$ws.Cells.Item(27, 4).Value = 82000
